# Update Bank Deposit data - append four new deposit rows (24-27) for
# collection date 05-12-2025, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = "05-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "Cash"; D = "2025-12-05"; E = 28180; F = ""; G = "2025-12-23" },
    @{ A = "05-12-2025"; B = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"; C = "Cash"; D = "2025-12-05"; E = 59730; F = ""; G = "2025-12-23" },
    @{ A = "05-12-2025"; B = "020965018-Kai Lalsingrao Shinde Gr.Big.Sheti Sah.Pat.Ltd. Br. Medha"; C = "Cash"; D = "2025-12-05"; E = 32740; F = ""; G = "2025-12-23" },
    @{ A = "05-12-2025"; B = "020965021-KAI.LALSINGRAO BAPUSO SHINDE SAH.PAT.LTD.,KUDAL, BR.KARAHAR"; C = "Cash"; D = "2025-12-05"; E = 4020; F = ""; G = "2025-12-23" }
)

$startRow = 24

# Columns A, D and G hold dates written as plain text (e.g. "05-12-2025",
# "2025-12-05") in this workbook, just like the pre-existing rows 2-23.
# Typing those strings straight into a General cell would make Excel's
# smart-input turn them into real date serials, so mark the destination
# cells as Text first, write the literal strings, then drop the explicit
# "Text" number-format again (matching the unstyled look of the original
# cells) once the value has been committed as text.
$textCols = @(1, 4, 7)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    foreach ($col in $textCols) {
        $ws.Cells.Item($r, $col).NumberFormat = "@"
    }
}

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    foreach ($col in $textCols) {
        $ws.Cells.Item($r, $col).Style = "Normal"
    }
}
